$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append the new row 33 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(33, 1).Value = "Demo inplannen"
$logs.Cells.Item(33, 2).Value = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item(33, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item(33, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item(33, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Cells.Item(33, 6).Value = "2025-08-14 21:55:45"
$logs.Cells.Item(33, 7).Value = "Nee"
$logs.Cells.Item(33, 8).Value = "Ja"
$logs.Cells.Item(33, 9).Value = "Nee"
$logs.Cells.Item(33, 10).Value = "Nee"

# --- "Logs" sheet: extend conditional-formatting ranges to include row 33 ---
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range("$col`2:$col`32")
    $newRange = $logs.Range("$col`2:$col`33")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- "Dashboard" sheet: bump the count for this category ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 25
